# Daily attendance processing - 2025-12-29 06:45:17
# Reorders the "Recorded By" names/emails in column G: for every row whose
# value contains multiple comma-separated entries, the entries are
# reversed in order (first becomes last, etc.). Single-entry cells are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        if ($value -like "*,*") {
            $parts = $value -split ","
            $trimmed = @()
            foreach ($p in $parts) {
                $trimmed += $p.Trim()
            }
            $reversed = $trimmed[($trimmed.Count - 1)..0]
            $newValue = [string]::Join(", ", $reversed)
            $cell.Value2 = $newValue
        }
    }
}
